$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for the two new columns, matching the style of the existing header row
$ws.Range("O1").Value = "Operacion"
$ws.Range("P1").Value = "Zona"
$headerRange = $ws.Range("O1:P1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Fix a typo in the address of case 7 (row 7)
$ws.Range("C7").Value = "Anchorena 1288"

# Populate Operacion (O) / Zona (P) for each data row
$ws.Range("O2").Value = "Recoleta"
$ws.Range("P2").Value = "Capital Sur"
$ws.Range("O3").Value = "Paternal"
$ws.Range("P3").Value = "Capital Norte"
$ws.Range("O4").Value = "Devoto"
$ws.Range("P4").Value = "Capital Norte"
$ws.Range("O5").Value = "Colegiales"
$ws.Range("P5").Value = "Capital Norte"
$ws.Range("O6").Value = "Devoto"
$ws.Range("P6").Value = "Capital Norte"
$ws.Range("O7").Value = "Almagro"
$ws.Range("P7").Value = "Capital Sur"
$ws.Range("O8").Value = "Recoleta"
$ws.Range("P8").Value = "Capital Sur"
$ws.Range("O9").Value = "Recoleta"
$ws.Range("P9").Value = "Capital Sur"
$ws.Range("O10").Value = "Paternal"
$ws.Range("P10").Value = "Capital Norte"
$ws.Range("O11").Value = "Saavedra"
$ws.Range("P11").Value = "Capital Norte"
$ws.Range("O12").Value = "Palermo"
$ws.Range("P12").Value = "Capital Sur"
$ws.Range("O13").Value = "San Telmo"
$ws.Range("P13").Value = "Capital Sur"
$ws.Range("O14").Value = "Palermo"
$ws.Range("P14").Value = "Capital Sur"
$ws.Range("O15").Value = "Colegiales"
$ws.Range("P15").Value = "Capital Norte"
$ws.Range("O16").Value = "Paternal"
$ws.Range("P16").Value = "Capital Norte"
$ws.Range("O17").Value = "Devoto"
$ws.Range("P17").Value = "Capital Norte"
$ws.Range("O18").Value = "Boedo"
$ws.Range("P18").Value = "Capital Sur"
$ws.Range("O19").Value = "Paternal"
$ws.Range("P19").Value = "Capital Norte"
$ws.Range("O20").Value = "Palermo"
$ws.Range("P20").Value = "Capital Sur"
$ws.Range("O21").Value = "Paternal"
$ws.Range("P21").Value = "Capital Norte"
$ws.Range("O22").Value = "Palermo"
$ws.Range("P22").Value = "Capital Sur"
$ws.Range("O23").Value = "Devoto"
$ws.Range("P23").Value = "Capital Norte"
$ws.Range("O24").Value = "Devoto"
$ws.Range("P24").Value = "Capital Norte"
$ws.Range("O25").Value = "Recoleta"
$ws.Range("P25").Value = "Capital Sur"
$ws.Range("O26").Value = "San Telmo"
$ws.Range("P26").Value = "Capital Sur"
$ws.Range("O27").Value = "Almagro"
$ws.Range("P27").Value = "Capital Sur"
$ws.Range("O28").Value = "Saavedra"
$ws.Range("P28").Value = "Capital Norte"
$ws.Range("O29").Value = "Saavedra"
$ws.Range("P29").Value = "Capital Norte"
$ws.Range("O30").Value = "Colegiales"
$ws.Range("P30").Value = "Capital Norte"
$ws.Range("O31").Value = "Saavedra"
$ws.Range("P31").Value = "Capital Norte"
$ws.Range("O32").Value = "San Telmo"
$ws.Range("P32").Value = "Capital Sur"
$ws.Range("O33").Value = "Boedo"
$ws.Range("P33").Value = "Capital Sur"
$ws.Range("O34").Value = "Almagro"
$ws.Range("P34").Value = "Capital Sur"
$ws.Range("O35").Value = "San Telmo"
$ws.Range("P35").Value = "Capital Sur"
$ws.Range("O36").Value = "Palermo"
$ws.Range("P36").Value = "Capital Sur"
$ws.Range("O37").Value = "Boedo"
$ws.Range("P37").Value = "Capital Sur"
$ws.Range("O38").Value = "San Telmo"
$ws.Range("P38").Value = "Capital Sur"
$ws.Range("O39").Value = "Palermo"
$ws.Range("P39").Value = "Capital Sur"
$ws.Range("O40").Value = "Palermo"
$ws.Range("P40").Value = "Capital Sur"
$ws.Range("O41").Value = "Palermo"
$ws.Range("P41").Value = "Capital Sur"
$ws.Range("O42").Value = "Palermo"
$ws.Range("P42").Value = "Capital Sur"
$ws.Range("O43").Value = "Boedo"
$ws.Range("P43").Value = "Capital Sur"
$ws.Range("O44").Value = "Saavedra"
$ws.Range("P44").Value = "Capital Norte"
$ws.Range("O45").Value = "Colegiales"
$ws.Range("P45").Value = "Capital Norte"
$ws.Range("O46").Value = "Palermo"
$ws.Range("P46").Value = "Capital Sur"
$ws.Range("O47").Value = "Colegiales"
$ws.Range("P47").Value = "Capital Norte"
$ws.Range("O48").Value = "Almagro"
$ws.Range("P48").Value = "Capital Sur"
$ws.Range("O49").Value = "Saavedra"
$ws.Range("P49").Value = "Capital Norte"
$ws.Range("O50").Value = "Boedo"
$ws.Range("P50").Value = "Capital Sur"
$ws.Range("O51").Value = "Saavedra"
$ws.Range("P51").Value = "Capital Norte"
$ws.Range("O52").Value = "Recoleta"
$ws.Range("P52").Value = "Capital Sur"
$ws.Range("O53").Value = "San Telmo"
$ws.Range("P53").Value = "Capital Sur"
$ws.Range("O54").Value = "San Telmo"
$ws.Range("P54").Value = "Capital Sur"
$ws.Range("O55").Value = "San Telmo"
$ws.Range("P55").Value = "Capital Sur"
$ws.Range("O56").Value = "San Telmo"
$ws.Range("P56").Value = "Capital Sur"
